# Refresh the cryptocurrency price / 1h-volume table to the values
# captured on Fri Jun 28 03:38:37 UTC 2024 (also re-sorts a few rows
# whose rank changed: PancakeSwap/Fetch.AI/ICP and Filecoin/OKB/Stacks/FDUSD).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.873.40"
$ws.Range("E2").Value = "  +1.31%  "
# Row 3
$ws.Range("D3").Value = "3.462.95"
$ws.Range("E3").Value = "  +2.38%  "
# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "
# Row 5
$ws.Range("D5").Value = "'583.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.56%  "
# Row 6
$ws.Range("D6").Value = "'147.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.50%  "
# Row 7
$ws.Range("D7").Value = "3.464.95"
$ws.Range("E7").Value = "  +2.46%  "
# Row 8
$ws.Range("E8").Value = "  +0.01%  "
# Row 9
$ws.Range("D9").Value = "'0.476"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.35%  "
# Row 10
$ws.Range("D10").Value = "'7.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.24%  "
# Row 11
$ws.Range("E11").Value = "  +3.50%  "
# Row 12
$ws.Range("E12").Value = "  +2.84%  "
# Row 13
$ws.Range("D13").Value = "4.046.68"
$ws.Range("E13").Value = "  +2.21%  "
# Row 14
$ws.Range("D14").Value = "'28.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.48%  "
# Row 16
$ws.Range("E16").Value = "  +1.67%  "
# Row 17
$ws.Range("D17").Value = "3.462.27"
$ws.Range("E17").Value = "  +2.31%  "
# Row 18
$ws.Range("D18").Value = "61.953.61"
$ws.Range("E18").Value = "  +1.19%  "
# Row 19
$ws.Range("D19").Value = "'6.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.86%  "
# Row 20
$ws.Range("E20").Value = "  +4.02%  "
# Row 21
$ws.Range("D21").Value = "'9.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.52%  "
# Row 22
$ws.Range("D22").Value = "'390.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.15%  "
# Row 23
$ws.Range("D23").Value = "'0.567"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.89%  "
# Row 24
$ws.Range("D24").Value = "'73.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.94%  "
# Row 25
$ws.Range("E25").Value = "  +0.21%  "
# Row 26
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.39%  "
# Row 27
$ws.Range("D27").Value = "'0.0000124"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.59%  "
# Row 28
$ws.Range("D28").Value = "3.599.80"
$ws.Range("E28").Value = "  +2.12%  "
# Row 29
$ws.Range("D29").Value = "'0.182"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.23%  "
# Row 30
$ws.Range("D30").Value = "'7.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.79%  "
# Row 31
$ws.Range("D31").Value = "'0.996"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.27%  "
# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'2.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.34%  "
# Row 33
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.06%  "
# Row 34
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'8.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.99%  "
# Row 35
$ws.Range("E35").Value = "  +0.04%  "
# Row 36
$ws.Range("D36").Value = "'24.25"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.49%  "
# Row 37
$ws.Range("D37").Value = "3.487.90"
$ws.Range("E37").Value = "  +2.74%  "
# Row 38
$ws.Range("D38").Value = "'7.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.87%  "
# Row 39
$ws.Range("D39").Value = "'1.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.19%  "
# Row 40
$ws.Range("D40").Value = "'5.17"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.46%  "
# Row 41
$ws.Range("D41").Value = "'166.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.44%  "
# Row 42
$ws.Range("E42").Value = "  +3.75%  "
# Row 43
$ws.Range("D43").Value = "'27.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.20%  "
# Row 44
$ws.Range("D44").Value = "'0.809"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.43%  "
# Row 45
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").Value = "'4.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.93%  "
# Row 46
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "'42.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.02%  "
# Row 47
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "'1.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.25%  "
# Row 48
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").Value = "'0.999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.10%  "
# Row 49
$ws.Range("E49").Value = "  -1.53%  "
# Row 50
$ws.Range("D50").Value = "2.574.29"
$ws.Range("E50").Value = "  +2.29%  "
# Row 51
$ws.Range("D51").Value = "'6.95"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.44%  "
